# Add "objective status" (size_preference) output column and update the
# lead flags to reflect the solver's feasibility/optimality result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new column E ("size_preference") for every invigilator row
# (rows 2-76) with the computed objective status string "s,m,l".
$ws.Range("E2:E76").Value = "s,m,l"

# Rows 25-76 (invig_id 23 onward) are no longer selected as "lead" -
# flip their lead flag from 1 to 0.
$ws.Range("D25:D76").Value = 0

# Move the active selection from F68 to F67, matching the saved view.
$ws.Range("F67").Select()
